$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Range("A52").Value = "x"
$rng = $ws.Range("A52")
$rng.Interior.ColorIndex = -4142
$rng.Borders.Color = 0
Write-Host "done"
